{"js": "// The document originally has a single run \"MiU 1303 \u2013 03/14\" (immediately\n// followed by the \"_GoBack\" bookmark, then a trailing \"/13\" run). Word's\n// proofer later re-checked that text and split it into four runs, wrapping\n// \"MiU\" and \"1303 \u2013 \" each in a w:proofErr spellStart/spellEnd pair, while\n// leaving \" \" and \"03/14\" as plain runs \u2014 the visible text is unchanged.\n\nconst body = context.document.body;\n\n// Locate the exact run text we need to restructure.\nconst matches = body.search(\"MiU 1303 \\u2013 03/14\", { matchCase: true, matchWholeWord: false });\nmatches.load(\"items\");\nawait context.sync();\n\nif (matches.items.length === 0) {\n  throw new Error(\"Could not find the target text 'MiU 1303 \\u2013 03/14'.\");\n}\n\nconst target = matches.items[0];\n\n// Flat-OPC OOXML fragment reproducing the run split + proofErr markers.\n// Word's InsertXML/insertOoxml always *replaces* the given range's content,\n// so the \"03/14\" tail has to be part of the payload too (it keeps its\n// original, plain-run shape).\nconst ooxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body><w:p>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r><w:t>MiU</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r><w:t xml:space=\"preserve\">1303 \\u2013 </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:t>03/14</w:t></w:r>' +\n  '</w:p></w:body></w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\n\n// Note whether the \"_GoBack\" bookmark exists (and is anchored right after\n// our target run) before we touch anything, so we can restore it precisely.\nconst bookmark = body.bookmarks.getByNameOrNullObject(\"_GoBack\");\nbookmark.load(\"isNullObject\");\nawait context.sync();\nconst hadGoBackBookmark = !bookmark.isNullObject;\n\ntarget.insertOoxml(ooxml, Word.InsertLocation.replace);\nawait context.sync();\n\nif (hadGoBackBookmark) {\n  // The replace above drags the pre-existing \"_GoBack\" bookmark to the\n  // front of the newly inserted content (it sat immediately after the old\n  // run, so it gets pushed to the start of the replacement). Put it back\n  // exactly where it was: right after \"03/14\", before the trailing \"/13\"\n  // run.\n  context.document.deleteBookmark(\"_GoBack\");\n  await context.sync();\n\n  const tail = body.search(\"03/14\", { matchCase: true });\n  await context.sync();\n\n  if (tail.items.length === 0) {\n    throw new Error(\"Could not find '03/14' to reattach the _GoBack bookmark.\");\n  }\n\n  const tailEnd = tail.items[0].getRange(\"End\");\n  tailEnd.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# The paragraph originally holds a single run \"MiU 1303 - 03/14\" (an en\n# dash, not a hyphen) immediately followed by the \"_GoBack\" bookmark and a\n# trailing \"/13\" run. Word's proofer later re-checked that text and split\n# it into four runs, wrapping \"MiU\" and \"1303 - \" each in a w:proofErr\n# spellStart/spellEnd pair, while \" \" and \"03/14\" stay plain runs. The\n# visible text itself does not change.\n\n$d = $word.ActiveDocument\n$enDash = [char]0x2013\n\n# Locate the exact text we need to restructure.\n$find = $d.Content\n$found = $find.Find.Execute(\"MiU 1303 \" + $enDash + \" 03/14\")\nif (-not $found) {\n    throw \"Could not find the target text 'MiU 1303 $enDash 03/14'.\"\n}\n\n# Re-materialize the match as a plain Range (the Find-mutated range object\n# inserts *after* itself instead of replacing, so build a fresh Range with\n# the same Start/End before calling InsertXML).\n$target = $d.Range($find.Start, $find.End)\n\n# Flat-OPC OOXML fragment reproducing the run split + proofErr markers.\n# InsertXML always replaces the given range's content, so the trailing\n# \"03/14\" has to be included in the payload too (kept as a plain run).\n$ooxml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body><w:p>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>MiU</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t xml:space=\"preserve\">1303 ' + $enDash + ' </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t>03/14</w:t></w:r>' +\n    '</w:p></w:body></w:document>' +\n    '</pkg:xmlData></pkg:part></pkg:package>'\n$target.InsertXML($ooxml)\n\n# The replace above drags the pre-existing \"_GoBack\" bookmark to the front\n# of the newly inserted content (it sat immediately after the old run, so\n# it gets pushed to the start of the replacement). Put it back exactly\n# where it was: right after \"03/14\", before the trailing \"/13\" run.\n$bm = $d.Bookmarks(\"_GoBack\")\n$bm.Delete()\n\n$tail = $d.Content\n$tailFound = $tail.Find.Execute(\"03/14\")\nif (-not $tailFound) {\n    throw \"Could not find '03/14' to reattach the _GoBack bookmark.\"\n}\n$insertPoint = $d.Range($tail.End, $tail.End)\n$d.Bookmarks.Add(\"_GoBack\", $insertPoint)\n"}
